$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 467
$ws.Range("F6").Value = 220
$ws.Range("F7").Value = 199
$ws.Range("F8").Value = 230
$ws.Range("F9").Value = 2816
$ws.Range("F11").Value = 115
$ws.Range("F12").Value = 2156
$ws.Range("F13").Value = 231
$ws.Range("F17").Value = 2505
$ws.Range("F19").Value = 1236
$ws.Range("F20").Value = 4505
$ws.Range("F22").Value = 4205
$ws.Range("F23").Value = 1252
$ws.Range("F24").Value = 2755
$ws.Range("F25").Value = 3166
$ws.Range("F26").Value = 127
$ws.Range("F27").Value = 1472
$ws.Range("F30").Value = 81
$ws.Range("F31").Value = 229
$ws.Range("F32").Value = 820
$ws.Range("F33").Value = 1394
$ws.Range("F34").Value = 106
$ws.Range("F35").Value = 220
$ws.Range("F36").Value = 553
$ws.Range("F38").Value = 269
$ws.Range("F39").Value = 339

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 87
$ws.Range("F8").Value = 83
$ws.Range("F10").Value = 17
$ws.Range("F15").Value = 87
$ws.Range("F16").Value = 49
$ws.Range("F17").Value = 45

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 467
$ws.Range("F5").Value = 87
$ws.Range("F8").Value = 220
$ws.Range("F9").Value = 199
$ws.Range("F11").Value = 230
$ws.Range("F12").Value = 2816
$ws.Range("F14").Value = 115
$ws.Range("F15").Value = 2157
$ws.Range("F16").Value = 231
$ws.Range("F20").Value = 17
$ws.Range("F21").Value = 2505
$ws.Range("F22").Value = 1236
$ws.Range("F26").Value = 4505
$ws.Range("F28").Value = 4205
$ws.Range("F29").Value = 1252
$ws.Range("F30").Value = 2755
$ws.Range("F31").Value = 3166
$ws.Range("F32").Value = 127
$ws.Range("F34").Value = 87
$ws.Range("F35").Value = 1472
$ws.Range("F36").Value = 49
$ws.Range("F39").Value = 81
$ws.Range("F40").Value = 229
$ws.Range("F41").Value = 820
$ws.Range("F42").Value = 45
$ws.Range("F43").Value = 1394
$ws.Range("F44").Value = 106
$ws.Range("F45").Value = 220
$ws.Range("F46").Value = 553
$ws.Range("F48").Value = 269
$ws.Range("F49").Value = 339
